$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last row (row 23) with the refreshed October 2025 stats.
$ws.Range("B23").Value = 6334
$ws.Range("C23").Value = 995
$ws.Range("D23").Value = 5899548
$ws.Range("E23").Value = 931.4095358383328
$ws.Range("F23").Value = 8.682223747426221
$ws.Range("G23").Value = 3.537981269510926
$ws.Range("H23").Value = 26.41221295533043
